# ADDED - BLog Pt2
# Updates the "social_diversidad_inclusion" workbook:
#  - Sheet 1 (DistPersonalxPto): fixes "Direccón" typo and refreshes the
#    headcount/percentage figures for Operación/Gerencias/Dirección.
#  - Sheet 2 (DistJuntaDirectivaXgenero): untouched.
#  - Sheet 3 (NuevasContratXgenero): replaces the single "11-100%" figure
#    with the new gender breakdown for new hires (two rows).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: DistPersonalxPto
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("DistPersonalxPto")

# Fix the "Direccón" -> "Dirección" typo in the header row.
$ws1.Range("C2").Value = "Dirección"

# Row 3 (Operación)
$ws1.Range("A3").Value = "131 – 53%"
$ws1.Range("B3").Value = "14 – 61%"
$ws1.Range("C3").Value = "2 – 40%"
$ws1.Range("D3").Value = 147

# Row 4 (Gerencias)
$ws1.Range("A4").Value = "118 – 47%"
$ws1.Range("B4").Value = "9 – 39%"
$ws1.Range("C4").Value = "3 – 60%"
$ws1.Range("D4").Value = 130

$ws1.Range("D3").Select() | Out-Null

# ---------------------------------------------------------------------
# Sheet 2: DistJuntaDirectivaXgenero (no content changes)
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("DistJuntaDirectivaXgenero")

# ---------------------------------------------------------------------
# Sheet 3: NuevasContratXgenero
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("NuevasContratXgenero")

$ws3.Range("A2").Value = "46 – 48%"
$ws3.Range("B2").ClearContents()
$ws3.Range("B2").Font.Bold = $false
$ws3.Range("A3").Value = "50 – 52%"

$ws3.Range("A4").Select() | Out-Null
